$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 994
$ws.Range("C3").Value = 988

$ws.Range("A4").Value = "Ron viejo de caldas (5años) botella"
$ws.Range("B4").Value = "Rones"
$ws.Range("C4").Value = 999
$ws.Range("D4").Value = 132000
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3/2/2026"
$ws.Range("E4").ClearFormats()
$ws.Range("F4").Value = 78000
